$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update generation Date value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-10-29T11:46:56+00:00"

# --- Elements sheet: fix casing of ExerciceProfessionnel references ---
$els = $wb.Worksheets.Item("Elements")

# Row 6 corresponds to CompetenceExclusive.exerciceProfessionnel
$els.Range("A6").Value = "CompetenceExclusive.ExerciceProfessionnel"
$els.Range("B6").Value = "CompetenceExclusive.ExerciceProfessionnel"
$els.Range("L6").Value = "Lien vers la classe ExerciceProfessionnel"
$els.Range("M6").Value = "Lien vers la classe ExerciceProfessionnel"
$els.Range("AF6").Value = "SavoirFaire.ExerciceProfessionnel"
